$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$preprocess = 'convert unicode to ascii, trim "space" and ",", convert to lower, remove multiple spaces, space after punctuation, remove break line'
$features   = '14 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), #digit/#ascii, %kwName, %kwAddress, %kwPhone, #max_digit_skip_0 >= 7, #max_digit_skip_0 = 0, #max_ascii_skip_0 >= 7, #max_ascii_skip_0 = 0, first_character_ascii, first_character_digit, last_character_ascii, last_character_digit'
$model      = 'Neuron Network'
$modelDetails = '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000'
$filter     = '0 filters: '

$rows = @(
    @{ Row=18; A='20160418_074755'; B=1201.699; G=0.991333333333333; H=0.940594059405941; J=0.0975609756097561 },
    @{ Row=19; A='20160418_080756'; B=1224.206;  G=0.996;             H=0.940594059405941; J=0.109756097560976 },
    @{ Row=20; A='20160418_082821'; B=1233.956;  G=0.990666666666667; H=0.95049504950495;  J=0.0470588235294118 },
    @{ Row=21; A='20160418_084854'; B=1253.958;  G=0.994;             H=0.95049504950495;  J=0.0588235294117647 },
    @{ Row=22; A='20160418_090948'; B=1254.064;  G=0.992;             H=0.940594059405941; J=0.0853658536585366 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $preprocess
    $ws.Cells.Item($row, 4).Value = $features
    $ws.Cells.Item($row, 5).Value = $model
    $ws.Cells.Item($row, 6).Value = $modelDetails
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $filter
    $ws.Cells.Item($row, 10).Value = $r.J
}
